$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "30.139.05"
$c.ClearFormats()
$ws.Range("E2").Value = "  -0.62%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.857.70"
$c.ClearFormats()
$ws.Range("E3").Value = "  -0.65%  "

$ws.Range("E4").Value = "  -0.05%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "232.81"
$c.ClearFormats()
$ws.Range("E5").Value = "  -0.96%  "

$ws.Range("E6").Value = "  +0.01%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4697"
$c.ClearFormats()
$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "42.68"
$c.ClearFormats()
$ws.Range("E8").Value = "  -0.52%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.2829"
$c.ClearFormats()
$ws.Range("E9").Value = "  -1.46%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.06446"
$c.ClearFormats()
$ws.Range("E10").Value = "  -1.90%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.73"
$c.ClearFormats()
$ws.Range("E11").Value = "  -4.11%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07662"
$c.ClearFormats()
$ws.Range("E12").Value = "  -4.42%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.848.92"
$c.ClearFormats()
$ws.Range("E13").Value = "  -1.11%  "

$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "93.22"
$c.ClearFormats()
$ws.Range("E14").Value = "  -3.86%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.028"
$c.ClearFormats()
$ws.Range("E15").Value = "  -1.64%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.6750"
$c.ClearFormats()
$ws.Range("E16").Value = "  -1.28%  "

$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "266.76"
$c.ClearFormats()
$ws.Range("E17").Value = "  -1.00%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "30.099.20"
$c.ClearFormats()
$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "13.26"
$c.ClearFormats()
$ws.Range("E19").Value = "  -5.49%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()
$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.000007487"
$c.ClearFormats()
$ws.Range("E21").Value = "  -1.78%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "2.100.49"
$c.ClearFormats()
$ws.Range("E22").Value = "  -0.41%  "

$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.ClearFormats()
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.126"
$c.ClearFormats()
$ws.Range("E24").Value = "  -2.93%  "

$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "6.063"
$c.ClearFormats()
$ws.Range("E25").Value = "  -2.50%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.234"
$c.ClearFormats()
$ws.Range("E26").Value = "  -2.03%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "165.21"
$c.ClearFormats()
$ws.Range("E27").Value = "  -2.15%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "18.45"
$c.ClearFormats()
$ws.Range("E28").Value = "  -2.40%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.872"
$c.ClearFormats()
$ws.Range("E29").Value = "  -3.86%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.367"
$c.ClearFormats()
$ws.Range("E30").Value = "  -0.22%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.09784"
$c.ClearFormats()
$ws.Range("E31").Value = "  -1.49%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.446"
$c.ClearFormats()
$ws.Range("E32").Value = "  -1.09%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.171"
$c.ClearFormats()
$ws.Range("E33").Value = "  -4.41%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.956"
$c.ClearFormats()
$ws.Range("E34").Value = "  -2.89%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.04628"
$c.ClearFormats()
$ws.Range("E35").Value = "  -1.51%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.106"
$c.ClearFormats()
$ws.Range("E36").Value = "  -2.84%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.6802"
$c.ClearFormats()
$ws.Range("E37").Value = "  -2.93%  "

$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.714"
$c.ClearFormats()
$ws.Range("E38").Value = "  +0.47%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01807"
$c.ClearFormats()
$ws.Range("E39").Value = "  -3.58%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.709"
$c.ClearFormats()
$ws.Range("E40").Value = "  +2.75%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "6.251"
$c.ClearFormats()
$ws.Range("E41").Value = "  -0.81%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "69.99"
$c.ClearFormats()
$ws.Range("E42").Value = "  -2.58%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()
$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.8275"
$c.ClearFormats()
$ws.Range("E44").Value = "  -1.65%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "102.11"
$c.ClearFormats()
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.863"
$c.ClearFormats()
$ws.Range("E46").Value = "  -5.29%  "

$ws.Range("B47").Value = "TheSandbox"
$ws.Range("C47").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.4024"
$c.ClearFormats()
$ws.Range("E47").Value = "  -3.43%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "9.160"
$c.ClearFormats()
$ws.Range("E48").Value = "  -0.52%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "918.60"
$c.ClearFormats()
$ws.Range("E49").Value = "  +0.57%  "

$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "6.865"
$c.ClearFormats()
$ws.Range("E50").Value = "  -2.66%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "34.02"
$c.ClearFormats()
$ws.Range("E51").Value = "  -1.31%  "
